# Added basic functional form/radio button for occupied and empty
# Adds a new column F (values default to 0 = "empty") for rows 1-21,
# matching the existing data range, and selects the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new column F with 0 for all data rows (1 through 21)
$ws.Range("F1:F21").Value = 0

# Reflect the new selection (whole column F of data) as the active view state
$ws.Range("F1:F21").Select()
